# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-14
$kValues = @{
    2  = 3
    3  = 3
    4  = 5
    5  = 3
    6  = 3
    7  = 1
    8  = 2
    9  = 3
    10 = 4
    11 = 5
    12 = 3
    13 = 3
    14 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
